$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 19453.416
$ws.Range("J69").Value = 20403.727
$ws.Range("L69").Value = 61211.181
$ws.Range("N69").Value = -62959.181

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 19453.416
$ws.Range("J72").Value = 20403.727
$ws.Range("L72").Value = 183633.543
$ws.Range("N72").Value = -192369.543

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4950
$ws.Range("J112").Value = 4950
$ws.Range("L112").Value = 14850
$ws.Range("N112").Value = -17066

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4950
$ws.Range("I116").Value = 4950
$ws.Range("K116").Value = 4950
$ws.Range("M116").Value = -1508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6378.2
$ws.Range("I137").Value = 5972.75
$ws.Range("K137").Value = 17918.25
$ws.Range("M137").Value = -15368.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3353.2683
$ws.Range("I32").Value = 3353.2683
$ws.Range("K32").Value = 3353.2683
$ws.Range("M32").Value = -3066.2683

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3230.889
$ws.Range("I61").Value = 3103.8667
$ws.Range("K61").Value = 3103.8667
$ws.Range("M61").Value = -2891.8667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1714.2858
$ws.Range("I74").Value = 1583.3334
$ws.Range("K74").Value = 1583.3334
$ws.Range("M74").Value = -709.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1714.2858
$ws.Range("I77").Value = 1583.3334
$ws.Range("K77").Value = 7916.666999999999
$ws.Range("M77").Value = -3548.666999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8450
$ws.Range("I122").Value = 8450
$ws.Range("K122").Value = 25350
$ws.Range("M122").Value = -22900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2668.6538
$ws.Range("I132").Value = 2655.4
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 7966.200000000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5436.200000000001
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3230.889
$ws.Range("I136").Value = 3103.8667
$ws.Range("K136").Value = 9311.6001
$ws.Range("M136").Value = -6761.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3033.7334
$ws.Range("I86").Value = 3013.3845
$ws.Range("K86").Value = 3013.3845
$ws.Range("M86").Value = -1890.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3033.7334
$ws.Range("I89").Value = 3013.3845
$ws.Range("K89").Value = 15066.9225
$ws.Range("M89").Value = -9450.922500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1307.3572
$ws.Range("I99").Value = 960.3333
$ws.Range("J99").Value = 1932
$ws.Range("K99").Value = 960.3333
$ws.Range("L99").Value = 1932
$ws.Range("M99").Value = 537.6667
$ws.Range("N99").Value = -4928

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7017.0835
$ws.Range("I134").Value = 7063.6313
$ws.Range("J134").Value = 6840.2
$ws.Range("K134").Value = 21190.8939
$ws.Range("L134").Value = 20520.6
$ws.Range("M134").Value = -18655.8939
$ws.Range("N134").Value = -25590.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.14285
$ws.Range("I7").Value = 73.5
$ws.Range("K7").Value = 73.5
$ws.Range("M7").Value = 39.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 4444893
$ws.Range("I22").Value = 253
$ws.Range("J22").Value = 6667213.5
$ws.Range("K22").Value = 253
$ws.Range("L22").Value = 6667213.5
$ws.Range("M22").Value = 97
$ws.Range("N22").Value = -6667913.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3033
$ws.Range("I58").Value = 2155.2222
$ws.Range("K58").Value = 2155.2222
$ws.Range("M58").Value = -1952.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1725.25
$ws.Range("J99").Value = 1603.5
$ws.Range("L99").Value = 1603.5
$ws.Range("N99").Value = -4599.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2267.4546
$ws.Range("I105").Value = 2586.8
$ws.Range("K105").Value = 2586.8
$ws.Range("M105").Value = -839.8000000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 383.46667
$ws.Range("I107").Value = 354.58334
$ws.Range("K107").Value = 354.58334
$ws.Range("M107").Value = 1565.41666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1725.25
$ws.Range("J126").Value = 1603.5
$ws.Range("L126").Value = 4810.5
$ws.Range("N126").Value = -9750.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1617.6765
$ws.Range("I134").Value = 1606.091
$ws.Range("K134").Value = 4818.272999999999
$ws.Range("M134").Value = -2283.272999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3033
$ws.Range("I136").Value = 2155.2222
$ws.Range("K136").Value = 6465.6666
$ws.Range("M136").Value = -3915.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 64.46154
$ws.Range("I12").Value = 56.8
$ws.Range("J12").Value = 69.25
$ws.Range("K12").Value = 170.4
$ws.Range("L12").Value = 207.75
$ws.Range("M12").Value = 2.600000000000023
$ws.Range("N12").Value = -553.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 625
$ws.Range("I17").Value = 175.33333
$ws.Range("J17").Value = 962.25
$ws.Range("K17").Value = 525.99999
$ws.Range("L17").Value = 2886.75
$ws.Range("M17").Value = -356.99999
$ws.Range("N17").Value = -3224.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 479.3
$ws.Range("J23").Value = 534.7143
$ws.Range("L23").Value = 1604.1429
$ws.Range("N23").Value = -2074.1429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1564.2941
$ws.Range("J26").Value = 998.8333
$ws.Range("L26").Value = 2996.4999
$ws.Range("N26").Value = -3572.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 4020
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 300
$ws.Range("M31").Value = -14712
$ws.Range("N31").Value = -876

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 3769
$ws.Range("J32").Value = 3930.7273
$ws.Range("L32").Value = 11792.1819
$ws.Range("N32").Value = -12358.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1731.6666
$ws.Range("J97").Value = 1997
$ws.Range("L97").Value = 5991
$ws.Range("N97").Value = -6983

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1920.7969
$ws.Range("I131").Value = 1499.5
$ws.Range("J131").Value = 1934.3871
$ws.Range("K131").Value = 4498.5
$ws.Range("L131").Value = 5803.1613
$ws.Range("M131").Value = 541.5
$ws.Range("N131").Value = -15883.1613

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8698.666999999999
$ws.Range("I70").Value = 8257.799999999999
$ws.Range("J70").Value = 9249.75
$ws.Range("K70").Value = 8257.799999999999
$ws.Range("L70").Value = 9249.75
$ws.Range("M70").Value = -7987.799999999999
$ws.Range("N70").Value = -9789.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8698.666999999999
$ws.Range("I73").Value = 8257.799999999999
$ws.Range("J73").Value = 9249.75
$ws.Range("K73").Value = 8257.799999999999
$ws.Range("L73").Value = 9249.75
$ws.Range("M73").Value = -7321.799999999999
$ws.Range("N73").Value = -11121.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2221
$ws.Range("I80").Value = 2221
$ws.Range("K80").Value = 2221
$ws.Range("M80").Value = -1223

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2221
$ws.Range("I83").Value = 2221
$ws.Range("K83").Value = 11105
$ws.Range("M83").Value = -6113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3391.875
$ws.Range("I126").Value = 2431.25
$ws.Range("J126").Value = 4352.5
$ws.Range("K126").Value = 7293.75
$ws.Range("L126").Value = 13057.5
$ws.Range("M126").Value = -4823.75
$ws.Range("N126").Value = -17997.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2680
$ws.Range("I132").Value = 2680
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8040
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5510
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3656.1333
$ws.Range("I7").Value = 3510.923
$ws.Range("J7").Value = 4600
$ws.Range("K7").Value = 3510.923
$ws.Range("L7").Value = 4600
$ws.Range("M7").Value = -3398.923
$ws.Range("N7").Value = -4824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 459.16666
$ws.Range("J16").Value = 448
$ws.Range("L16").Value = 448
$ws.Range("N16").Value = -788

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1104.8182
$ws.Range("I22").Value = 885.0769
$ws.Range("J22").Value = 1422.2222
$ws.Range("K22").Value = 885.0769
$ws.Range("L22").Value = 1422.2222
$ws.Range("M22").Value = -590.0769
$ws.Range("N22").Value = -2012.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1104.8182
$ws.Range("I27").Value = 885.0769
$ws.Range("J27").Value = 1422.2222
$ws.Range("K27").Value = 885.0769
$ws.Range("L27").Value = 1422.2222
$ws.Range("M27").Value = -778.0769
$ws.Range("N27").Value = -1636.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2909.6
$ws.Range("I68").Value = 2749
$ws.Range("K68").Value = 2749
$ws.Range("M68").Value = -2000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2909.6
$ws.Range("I71").Value = 2749
$ws.Range("K71").Value = 13745
$ws.Range("M71").Value = -10001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1994.5555
$ws.Range("I82").Value = 1524.5
$ws.Range("J82").Value = 2370.6
$ws.Range("K82").Value = 1524.5
$ws.Range("L82").Value = 2370.6
$ws.Range("M82").Value = -1163.5
$ws.Range("N82").Value = -3092.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1994.5555
$ws.Range("I85").Value = 1524.5
$ws.Range("J85").Value = 2370.6
$ws.Range("K85").Value = 1524.5
$ws.Range("L85").Value = 2370.6
$ws.Range("M85").Value = -276.5
$ws.Range("N85").Value = -4866.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2170.8572
$ws.Range("I93").Value = 2216.1667
$ws.Range("J93").Value = 1899
$ws.Range("K93").Value = 2216.1667
$ws.Range("L93").Value = 1899
$ws.Range("M93").Value = -968.1667000000002
$ws.Range("N93").Value = -4395

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1273.7273
$ws.Range("I100").Value = 1066
$ws.Range("J100").Value = 2208.5
$ws.Range("K100").Value = 1066
$ws.Range("L100").Value = 2208.5
$ws.Range("M100").Value = -525
$ws.Range("N100").Value = -3290.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3656.1333
$ws.Range("I126").Value = 3510.923
$ws.Range("J126").Value = 4600
$ws.Range("K126").Value = 10532.769
$ws.Range("L126").Value = 13800
$ws.Range("M126").Value = -8062.769
$ws.Range("N126").Value = -18740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1822.7097
$ws.Range("I132").Value = 1713.619
$ws.Range("K132").Value = 5140.857
$ws.Range("M132").Value = -2610.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4519.706
$ws.Range("I126").Value = 3910.7693
$ws.Range("J126").Value = 6498.75
$ws.Range("K126").Value = 11732.3079
$ws.Range("L126").Value = 19496.25
$ws.Range("M126").Value = -9262.3079
$ws.Range("N126").Value = -24436.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6174.067
$ws.Range("I132").Value = 5257.5713
$ws.Range("K132").Value = 15772.7139
$ws.Range("M132").Value = -13242.7139
